$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: "<name>_old" -> "<name>_FV2210", "<name>_new" -> "<name>_FV2304"
#    Columns A-J (1-10) carry the "_old" headers, column K (11) is "diff" (unchanged),
#    columns L-U (12-21) carry the "_new" headers.
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2304"
}

# 2) Turn the whole used range into an Excel table ("Table1") with headers.
$tableRange = $ws.Range("A1:U87")
$tbl = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split after row 1).
$ws.Activate()
[void]$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true

Write-Host "done"
